$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: 2014/12 (IFRS연결) ---
$ws.Cells.Item(2, 4).Value = 1363
$ws.Cells.Item(2, 5).Value = 205
$ws.Cells.Item(2, 6).Value = 205
$ws.Cells.Item(2, 7).Value = 176
$ws.Cells.Item(2, 8).Value = 117
$ws.Cells.Item(2, 9).Value = 117
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 2241
$ws.Cells.Item(2, 12).Value = 1003
$ws.Cells.Item(2, 13).Value = 1238
$ws.Cells.Item(2, 14).Value = 1235
$ws.Cells.Item(2, 15).Value = 3
$ws.Cells.Item(2, 16).Value = 133
$ws.Cells.Item(2, 17).Value = 203
$ws.Cells.Item(2, 18).Value = -73
$ws.Cells.Item(2, 19).Value = -53
$ws.Cells.Item(2, 20).Value = 99
$ws.Cells.Item(2, 21).Value = 104
$ws.Cells.Item(2, 22).Value = 653
$ws.Cells.Item(2, 23).Value = 15.05
$ws.Cells.Item(2, 24).Value = 8.57
$ws.Cells.Item(2, 25).Value = 9.84
$ws.Cells.Item(2, 26).Value = 5.33
$ws.Cells.Item(2, 27).Value = 81
$ws.Cells.Item(2, 28).Value = 837.48
$ws.Cells.Item(2, 29).Value = 3379
$ws.Cells.Item(2, 30).Value = 7.83
$ws.Cells.Item(2, 31).Value = 35630
$ws.Cells.Item(2, 32).Value = 0.74
$ws.Cells.Item(2, 33).Value = 956
$ws.Cells.Item(2, 34).Value = 3.61
$ws.Cells.Item(2, 35).Value = 28.28
$ws.Cells.Item(2, 36).Value = 3467892

# --- Row 3: 2015/12 (IFRS연결) ---
$ws.Cells.Item(3, 4).Value = 1358
$ws.Cells.Item(3, 5).Value = 148
$ws.Cells.Item(3, 6).Value = 148
$ws.Cells.Item(3, 7).Value = 208
$ws.Cells.Item(3, 8).Value = 158
$ws.Cells.Item(3, 9).Value = 158
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 2404
$ws.Cells.Item(3, 12).Value = 1090
$ws.Cells.Item(3, 13).Value = 1314
$ws.Cells.Item(3, 14).Value = 1309
$ws.Cells.Item(3, 15).Value = 5
$ws.Cells.Item(3, 16).Value = 133
$ws.Cells.Item(3, 17).Value = 189
$ws.Cells.Item(3, 18).Value = -103
$ws.Cells.Item(3, 19).Value = -44
$ws.Cells.Item(3, 20).Value = 134
$ws.Cells.Item(3, 21).Value = 55
$ws.Cells.Item(3, 22).Value = 664
$ws.Cells.Item(3, 23).Value = 10.87
$ws.Cells.Item(3, 24).Value = 11.63
$ws.Cells.Item(3, 25).Value = 12.43
$ws.Cells.Item(3, 26).Value = 6.8
$ws.Cells.Item(3, 27).Value = 82.92
$ws.Cells.Item(3, 28).Value = 927.79
$ws.Cells.Item(3, 29).Value = 4559
$ws.Cells.Item(3, 30).Value = 6.02
$ws.Cells.Item(3, 31).Value = 37758
$ws.Cells.Item(3, 32).Value = 0.73
$ws.Cells.Item(3, 33).Value = 1338
$ws.Cells.Item(3, 34).Value = 4.88
$ws.Cells.Item(3, 35).Value = 29.35
$ws.Cells.Item(3, 36).Value = 3467892

# --- Row 4: 2016/12 (IFRS연결) ---
$ws.Cells.Item(4, 4).Value = 1555
$ws.Cells.Item(4, 5).Value = 178
$ws.Cells.Item(4, 6).Value = 178
$ws.Cells.Item(4, 7).Value = 118
$ws.Cells.Item(4, 8).Value = 70
$ws.Cells.Item(4, 9).Value = 70
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 2379
$ws.Cells.Item(4, 12).Value = 1114
$ws.Cells.Item(4, 13).Value = 1265
$ws.Cells.Item(4, 14).Value = 1260
$ws.Cells.Item(4, 15).Value = 5
$ws.Cells.Item(4, 16).Value = 133
$ws.Cells.Item(4, 17).Value = 148
$ws.Cells.Item(4, 18).Value = -79
$ws.Cells.Item(4, 19).Value = -93
$ws.Cells.Item(4, 20).Value = 84
$ws.Cells.Item(4, 21).Value = 63
$ws.Cells.Item(4, 22).Value = 689
$ws.Cells.Item(4, 23).Value = 11.42
$ws.Cells.Item(4, 24).Value = 4.49
$ws.Cells.Item(4, 25).Value = 5.46
$ws.Cells.Item(4, 26).Value = 2.92
$ws.Cells.Item(4, 27).Value = 88.09
$ws.Cells.Item(4, 28).Value = 938.33
$ws.Cells.Item(4, 29).Value = 2023
$ws.Cells.Item(4, 30).Value = 14.62
$ws.Cells.Item(4, 31).Value = 37880
$ws.Cells.Item(4, 32).Value = 0.78
$ws.Cells.Item(4, 33).Value = 153
$ws.Cells.Item(4, 34).Value = 0.52
$ws.Cells.Item(4, 35).Value = 7.16
$ws.Cells.Item(4, 36).Value = 3467892

# --- Row 5: 2017/12 (IFRS연결) ---
$ws.Cells.Item(5, 4).Value = 1419
$ws.Cells.Item(5, 5).Value = 59
$ws.Cells.Item(5, 6).Value = 59
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(5, 8).Value = -18
$ws.Cells.Item(5, 9).Value = -17
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 2117
$ws.Cells.Item(5, 12).Value = 884
$ws.Cells.Item(5, 13).Value = 1232
$ws.Cells.Item(5, 14).Value = 1228
$ws.Cells.Item(5, 15).Value = 4
$ws.Cells.Item(5, 16).Value = 138
$ws.Cells.Item(5, 17).Value = 153
$ws.Cells.Item(5, 18).Value = -5
$ws.Cells.Item(5, 19).Value = -228
$ws.Cells.Item(5, 20).Value = 33
$ws.Cells.Item(5, 21).Value = 120
$ws.Cells.Item(5, 22).Value = 546
$ws.Cells.Item(5, 23).Value = 4.13
$ws.Cells.Item(5, 24).Value = -1.25
$ws.Cells.Item(5, 25).Value = -1.41
$ws.Cells.Item(5, 26).Value = -0.79
$ws.Cells.Item(5, 27).Value = 71.76000000000001
$ws.Cells.Item(5, 28).Value = 833.4400000000001
$ws.Cells.Item(5, 29).Value = -521
$ws.Cells.Item(5, 30).Value = -53.81
$ws.Cells.Item(5, 31).Value = 39201
$ws.Cells.Item(5, 32).Value = 0.72
$ws.Cells.Item(5, 33).Value = 635
$ws.Cells.Item(5, 34).Value = 2.26
$ws.Cells.Item(5, 35).Value = -113.88
$ws.Cells.Item(5, 36).Value = 3265735

# --- Row 6: 2018/12 (IFRS연결) ---
$ws.Cells.Item(6, 4).Value = 1323
$ws.Cells.Item(6, 5).Value = 30
$ws.Cells.Item(6, 6).Value = 30
$ws.Cells.Item(6, 7).Value = 28
$ws.Cells.Item(6, 8).Value = 14
$ws.Cells.Item(6, 9).Value = 14
$ws.Cells.Item(6, 11).Value = 1959
$ws.Cells.Item(6, 12).Value = 833
$ws.Cells.Item(6, 13).Value = 1127
$ws.Cells.Item(6, 14).Value = 1123
$ws.Cells.Item(6, 16).Value = 170
$ws.Cells.Item(6, 17).Value = 79
$ws.Cells.Item(6, 18).Value = 132
$ws.Cells.Item(6, 19).Value = -218
$ws.Cells.Item(6, 20).Value = 31
$ws.Cells.Item(6, 21).Value = 47
$ws.Cells.Item(6, 22).Value = 458
$ws.Cells.Item(6, 23).Value = 2.27
$ws.Cells.Item(6, 24).Value = 1.03
$ws.Cells.Item(6, 25).Value = 1.17
$ws.Cells.Item(6, 26).Value = 0.67
$ws.Cells.Item(6, 27).Value = 73.90000000000001
$ws.Cells.Item(6, 28).Value = 588.83
$ws.Cells.Item(6, 29).Value = 446
$ws.Cells.Item(6, 30).Value = 64.56999999999999
$ws.Cells.Item(6, 31).Value = 40725
$ws.Cells.Item(6, 32).Value = 0.71
$ws.Cells.Item(6, 33).Value = 1000
$ws.Cells.Item(6, 34).Value = 3.47
$ws.Cells.Item(6, 35).Value = 208.67
$ws.Cells.Item(6, 36).Value = 2895569

# --- Rows 7-9: 2019/12(E), 2020/12(E), 2021/12(E) estimates removed ---
$ws.Range("D7:AJ9").ClearContents()
